$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-13) but keep header row 1 intact (style, values).
$ws.Range("A2:E13").Clear()

# Populate columns A-D (text columns) column-by-column so that shared-string
# insertion order matches a column-major fill (brand, model, year, code, price).
# Column A
$ws.Range("A2").Value = "Acura"
$ws.Range("A3").Value = "Acura"
$ws.Range("A4").Value = "Acura"
$ws.Range("A5").Value = "Acura"
$ws.Range("A6").Value = "Acura"
$ws.Range("A7").Value = "Acura"
$ws.Range("A8").Value = "Acura"
$ws.Range("A9").Value = "Acura"
$ws.Range("A10").Value = "Acura"
$ws.Range("A11").Value = "Acura"
$ws.Range("A12").Value = "Acura"
$ws.Range("A13").Value = "Acura"
$ws.Range("A14").Value = "Acura"
$ws.Range("A15").Value = "Acura"
$ws.Range("A16").Value = "Acura"
$ws.Range("A17").Value = "Agrale"
$ws.Range("A18").Value = "Agrale"
$ws.Range("A19").Value = "Agrale"
$ws.Range("A20").Value = "Agrale"

# Column B
$ws.Range("B2").Value = "Integra GS 1.8"
$ws.Range("B3").Value = "Integra GS 1.8"
$ws.Range("B4").Value = "Legend 3.2/3.5"
$ws.Range("B5").Value = "Legend 3.2/3.5"
$ws.Range("B6").Value = "Legend 3.2/3.5"
$ws.Range("B7").Value = "Legend 3.2/3.5"
$ws.Range("B8").Value = "Legend 3.2/3.5"
$ws.Range("B9").Value = "Legend 3.2/3.5"
$ws.Range("B10").Value = "Legend 3.2/3.5"
$ws.Range("B11").Value = "Legend 3.2/3.5"
$ws.Range("B12").Value = "NSX 3.0"
$ws.Range("B13").Value = "NSX 3.0"
$ws.Range("B14").Value = "NSX 3.0"
$ws.Range("B15").Value = "NSX 3.0"
$ws.Range("B16").Value = "NSX 3.0"
$ws.Range("B17").Value = "MARRUÁ 2.8 12V 132cv TDI Diesel"
$ws.Range("B18").Value = "MARRUÁ 2.8 12V 132cv TDI Diesel"
$ws.Range("B19").Value = "MARRUÁ 2.8 12V 132cv TDI Diesel"
$ws.Range("B20").Value = "MARRUÁ 2.8 12V 132cv TDI Diesel"

# Column C
$ws.Range("C2").Value = "1992 Gasolina"
$ws.Range("C3").Value = "1991 Gasolina"
$ws.Range("C4").Value = "1998 Gasolina"
$ws.Range("C5").Value = "1997 Gasolina"
$ws.Range("C6").Value = "1996 Gasolina"
$ws.Range("C7").Value = "1995 Gasolina"
$ws.Range("C8").Value = "1994 Gasolina"
$ws.Range("C9").Value = "1993 Gasolina"
$ws.Range("C10").Value = "1992 Gasolina"
$ws.Range("C11").Value = "1991 Gasolina"
$ws.Range("C12").Value = "1995 Gasolina"
$ws.Range("C13").Value = "1994 Gasolina"
$ws.Range("C14").Value = "1993 Gasolina"
$ws.Range("C15").Value = "1992 Gasolina"
$ws.Range("C16").Value = "1991 Gasolina"
$ws.Range("C17").Value = "2007 Diesel"
$ws.Range("C18").Value = "2006 Diesel"
$ws.Range("C19").Value = "2005 Diesel"
$ws.Range("C20").Value = "2004 Diesel"

# Column D
$ws.Range("D2").Value = "038003-2"
$ws.Range("D3").Value = "038003-2"
$ws.Range("D4").Value = "038002-4"
$ws.Range("D5").Value = "038002-4"
$ws.Range("D6").Value = "038002-4"
$ws.Range("D7").Value = "038002-4"
$ws.Range("D8").Value = "038002-4"
$ws.Range("D9").Value = "038002-4"
$ws.Range("D10").Value = "038002-4"
$ws.Range("D11").Value = "038002-4"
$ws.Range("D12").Value = "038001-6"
$ws.Range("D13").Value = "038001-6"
$ws.Range("D14").Value = "038001-6"
$ws.Range("D15").Value = "038001-6"
$ws.Range("D16").Value = "038001-6"
$ws.Range("D17").Value = "060001-6"
$ws.Range("D18").Value = "060001-6"
$ws.Range("D19").Value = "060001-6"
$ws.Range("D20").Value = "060001-6"

# Column E (price) values look numeric ("  11120.00"), so force them to stay text
# (matching the shared-string <t> entries) by setting a Text number format first.
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = " 11120.00"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = " 10387.00"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = " 14248.00"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = " 33464.00"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = " 33464.00"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = " 33464.00"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = " 33464.00"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = " 33464.00"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = " 48348.00"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = " 44601.00"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = " 43449.00"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = " 36830.00"
